$d = $word.ActiveDocument

# --- Step 1: remove the whole "Histórico" bullet paragraph -----------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Histórico") {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: "Carla - Erros" -> "Carla - Erros 2 primeiros" ----------------
$d.Content.Find.Execute("Carla - Erros", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Carla - Erros 2 primeiros", 2)

# --- Step 3: relocate the _GoBack bookmark ----------------------------------
# It currently sits at the end of the "Unidades - DInheiro" paragraph; it
# should end up at the end of the "Carla - Erros 2 primeiros" paragraph
# (which is the paragraph right before it).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Carla - Erros 2 primeiros") {
        $target = $p.Range
        break
    }
}

$endPos = $target.End - 1

# NOTE: inserting a zero-length bookmark exactly at "paragraph end - 1" (i.e.
# right before the paragraph mark) mis-anchors it to the start of the
# document. Work around this by temporarily inserting a placeholder
# character after the intended bookmark position, adding the bookmark before
# it, and then removing the placeholder again.
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Delete()
